$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Settings sheet
# ---------------------------------------------------------------------
$settings = $wb.Worksheets.Item("Settings")

# Row 2: OrchestratorQueueName value -> new queue name
$settings.Cells.Item(2, 2).Value = "RM_ServiceRequestsReportingQueue"

# Row 5: logF_BusinessProcessName value -> new queue name (reuse same text)
$settings.Cells.Item(5, 2).Value = "RM_ServiceRequestsReportingQueue"

# New row 7: strSheetNameTransactionDetails / Transaction Details
$settings.Cells.Item(7, 1).Value = "strSheetNameTransactionDetails"
$settings.Cells.Item(7, 2).Value = "Transaction Details"

# New row 8: strSheeNameSummary / Transaction Outcomes Summary
$settings.Cells.Item(8, 1).Value = "strSheeNameSummary"
$settings.Cells.Item(8, 2).Value = "Transaction Outcomes Summary"

# ---------------------------------------------------------------------
# Assets sheet
# ---------------------------------------------------------------------
$assets = $wb.Worksheets.Item("Assets")

$assets.Cells.Item(2, 1).Value = "strDailyReportPath"
$assets.Cells.Item(2, 2).Value = "RM_Report_Path"

$assets.Cells.Item(3, 1).Value = "strDailyReportFileName"
$assets.Cells.Item(3, 2).Value = "RM_Report_FileName"

$assets.Cells.Item(4, 1).Value = "strDailyReportTemplatePath"
$assets.Cells.Item(4, 2).Value = "RM_Report_TemplatePath"

$assets.Cells.Item(5, 1).Value = "strDailyReportTemplateFileName"
$assets.Cells.Item(5, 2).Value = "RM_Report_TemplateFileName"
$assets.Cells.Item(5, 4).Value = "ReturnedMail_DailyReport_{0}"

# ---------------------------------------------------------------------
# View / selection state
# ---------------------------------------------------------------------
[void]$assets.Activate()
[void]$assets.Range("D2").Select()

[void]$settings.Activate()
[void]$settings.Range("C9").Select()
